# Update statistics for the "GO" state rows (DF add update per commit message)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: year 2015, state GO
$ws.Range("C2").Value = 12548
$ws.Range("E2").Value = 398.9145028330469
$ws.Range("F2").Value = 398.9312744140625
$ws.Range("H2").Value = 2.132369496561104
$ws.Range("I2").Value = 1.460263502440948
$ws.Range("J2").Value = 0.366059266351641
$ws.Range("K2").Value = -0.1592442942796279
$ws.Range("L2").Value = 0.7963502066840924

# Row 7: year 2016, state GO
$ws.Range("C7").Value = 14714
$ws.Range("E7").Value = 402.2692439815702
$ws.Range("F7").Value = 402.34423828125
$ws.Range("H7").Value = 1.935110439601139
$ws.Range("I7").Value = 1.391082470452827
$ws.Range("J7").Value = 0.3458088062324144
$ws.Range("K7").Value = -0.1066689038923072
$ws.Range("L7").Value = 0.7729680457828501

# Row 12: year 2017, state GO
$ws.Range("C12").Value = 8879
$ws.Range("E12").Value = 404.0198314970497
$ws.Range("F12").Value = 404.0214538574219
$ws.Range("H12").Value = 1.645235531872815
$ws.Range("I12").Value = 1.28266735043534
$ws.Range("J12").Value = 0.3174763341894781
$ws.Range("K12").Value = 0.2639047614576369
$ws.Range("L12").Value = 1.982068610517871

# Row 17: year 2018, state GO
$ws.Range("C17").Value = 12952
$ws.Range("E17").Value = 406.1755092778715
$ws.Range("F17").Value = 406.2969665527344
$ws.Range("H17").Value = 1.87967308393604
$ws.Range("I17").Value = 1.371011700874956
$ws.Range("J17").Value = 0.3375416955375869
$ws.Range("K17").Value = -0.5283813453827979
$ws.Range("L17").Value = 1.614059594151295

# Row 22: year 2019, state GO
$ws.Range("C22").Value = 11189
$ws.Range("E22").Value = 408.8485073219495
$ws.Range("F22").Value = 408.9463195800781
$ws.Range("H22").Value = 2.421529258660425
$ws.Range("I22").Value = 1.556126363333141
$ws.Range("J22").Value = 0.3806119712961951
$ws.Range("K22").Value = -0.3847638769677218
$ws.Range("L22").Value = 1.605163876050673

# Row 27: year 2020, state GO
$ws.Range("C27").Value = 13658
$ws.Range("E27").Value = 410.9726019404192
$ws.Range("F27").Value = 410.9277038574219
$ws.Range("H27").Value = 1.861093948299786
$ws.Range("I27").Value = 1.36421917165087
$ws.Range("J27").Value = 0.3319489341162085
$ws.Range("K27").Value = 0.08328362086594324
$ws.Range("L27").Value = 0.4182223542257963
